$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.007.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.01%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.741.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.60%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5784"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.75%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -1.45%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.17"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06585"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07531"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.75%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.743.91"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.715"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.39%  "

$ws.Range("E14").Value = "  -3.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.979.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "74.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.82%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008668"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -11.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.999.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.89%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.327"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.15%  "

$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "205.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.28%  "

$ws.Range("E23").Value = "  -3.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.066"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.83%  "

$ws.Range("E27").Value = "  -4.46%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.07%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.381"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.22%  "

$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06152"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.39%  "

$ws.Range("E31").Value = "  -3.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.743"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.03%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.713"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.671"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.03%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.035"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6372"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.416"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.71%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.717"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.97%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01668"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.128.26"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.30%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.200"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8733"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.62%  "

$ws.Range("E43").Value = "  +0.11%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.71%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.891.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.93%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000108"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.581"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.71%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.238"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05379"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.13%  "

$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.274"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.89%  "
